$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.2020408164059704
$ws.Range("F3").Value = 78.19999694824219
$ws.Range("G3").Value = 1563.999938964844
$ws.Range("H3").Value = 0.7979591835940296
$ws.Range("H6").Value = 0.3264929204294104
$ws.Range("F7").Value = 82.05000305175781
$ws.Range("G7").Value = 1641.000061035156
$ws.Range("H7").Value = 0.6735070795705895
$ws.Range("H8").Value = 0.1994233483076236
$ws.Range("F9").Value = 83.30000305175781
$ws.Range("G9").Value = 3332.000122070312
$ws.Range("H9").Value = 0.8005766516923764
$ws.Range("H10").Value = 0.2939456448799682
$ws.Range("F11").Value = 82.09999847412109
$ws.Range("G11").Value = 4104.999923706055
$ws.Range("H11").Value = 0.7060543551200318
$ws.Range("H13").Value = 0.2850789511996523
$ws.Range("F14").Value = 96.44999694824219
$ws.Range("G14").Value = 4822.499847412109
$ws.Range("H14").Value = 0.7149210488003477
$ws.Range("H15").Value = 0.1897198552766218
$ws.Range("F16").Value = 108.4499969482422
$ws.Range("G16").Value = 5422.499847412109
$ws.Range("H16").Value = 0.516962756928909
$ws.Range("H17").Value = 0.2933173877944692
$ws.Range("H18").Value = 0.1963294899032366
$ws.Range("F19").Value = 101.1999969482422
$ws.Range("G19").Value = 6071.999816894531
$ws.Range("H19").Value = 0.471190775767768
$ws.Range("H20").Value = 0.2392193399940886
$ws.Range("H21").Value = 0.09326039433490677
$ws.Range("H22").Value = 0.2342788964336497
$ws.Range("F23").Value = 97.55000305175781
$ws.Range("G23").Value = 5853.000183105469
$ws.Range("H23").Value = 0.4123393302057723
$ws.Range("H24").Value = 0.2655245458626109
$ws.Range("H25").Value = 0.08785722749796712
$ws.Range("H26").Value = 0.2200738230840038
$ws.Range("F27").Value = 92.44999694824219
$ws.Range("G27").Value = 5546.999816894531
$ws.Range("H27").Value = 0.3831605434493547
$ws.Range("H28").Value = 0.3101728691311063
$ws.Range("H29").Value = 0.0865927643355353
$ws.Range("H30").Value = 0.2105908844293644
$ws.Range("F31").Value = 97.75
$ws.Range("G31").Value = 5865
$ws.Range("H31").Value = 0.343445404846221
$ws.Range("H32").Value = 0.3732282411635292
$ws.Range("H33").Value = 0.0727354695608855
$ws.Range("H37").Value = 0.2756603970605329
$ws.Range("H38").Value = 0.5537965450563179
$ws.Range("H39").Value = 0.1134045900767662
$ws.Range("F40").Value = 137.5500030517578
$ws.Range("G40").Value = 687.7500152587891
$ws.Range("H40").Value = 0.05713846780638292
$ws.Range("H41").Value = 0.2675888384137958
$ws.Range("H42").Value = 0.4764331918316316
$ws.Range("H43").Value = 0.09540311040548617
$ws.Range("F44").Value = 148.8999938964844
$ws.Range("G44").Value = 2233.499908447266
$ws.Range("H44").Value = 0.1605748593490865
$ws.Range("H45").Value = 0.1957773201390204
$ws.Range("H46").Value = 0.5773148512662105
$ws.Range("H47").Value = 0.07288957832364869
$ws.Range("F48").Value = 136.6499938964844
$ws.Range("G48").Value = 2732.999877929688
$ws.Range("H48").Value = 0.1540182502711204
$ws.Range("H49").Value = 0.1605961443530085
$ws.Range("H50").Value = 0.5699450050133548
$ws.Range("H51").Value = 0.06105101779666405
$ws.Range("F52").Value = 130.5500030517578
$ws.Range("G52").Value = 3916.500091552734
$ws.Range("H52").Value = 0.2084078328369726
$ws.Range("H53").Value = 0.1182164291121852
$ws.Range("H54").Value = 0.6798089974915936
$ws.Range("H55").Value = 0.05301026098984984
$ws.Range("F56").Value = 76.94999694824219
$ws.Range("G56").Value = 2308.499908447266
$ws.Range("H56").Value = 0.1489643124063713
